# Applies the "Elimina EC anteriores y se agregan nuevos, se modifica base de datos" edit:
#  - Adds a new period row (1912) to the account-statement table, between the
#    existing first data row and the former second data row.
#  - Changes the first data row's period from 1911 to 2110 and updates its amounts.
#  - Updates the second data row's period/amount (now showing 1911).
#  - Updates the "VALOR MORA" total and "Cant. Periodos" count at the top of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 17 (pushes current row 17 -> 18, row 22 -> 23, row 23 -> 24, etc.)
$ws.Rows("17").Insert()

# Give the new row 17 the same formatting as row 16 (the first data row of the table)
$ws.Range("B16:J16").Copy()
$ws.Range("B17:J17").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Populate the new row 17 (period 1912) with the same "Tipo/No Doc/Nombre" as row 16
$ws.Range("B17").Value2 = $ws.Range("B16").Value2
$ws.Range("C17").Value2 = $ws.Range("C16").Value2
$ws.Range("D17").Value2 = $ws.Range("D16").Value2
$ws.Range("E17").Value2 = "1912"
$ws.Range("F17").Value2 = 120000
$ws.Range("G17").Value2 = 2000000

# Row 16 now represents period 2110, with updated amounts
$ws.Range("E16").Value2 = "2110"
$ws.Range("F16").Value2 = 4160
$ws.Range("G16").Value2 = 2000000

# Row 18 (the former row 17) now shows period 1911, with an updated "Salario Basico"
$ws.Range("E18").Value2 = "1911"
$ws.Range("G18").Value2 = 2000000

# Update the "VALOR MORA" total and "Cant. Periodos" count
$ws.Range("E11").Value2 = 244160
$ws.Range("F13").Value2 = 3
